$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.384999999999991
$ws.Range("C6").Value = -11.5449
$ws.Range("C7").Value = -12.0474
$ws.Range("B8").Value = 4.905799999999997
$ws.Range("C8").Value = -11.3917
$ws.Range("D11").Value = -9.037499999999993
$ws.Range("A12").Value = -22.77720000000001
$ws.Range("B12").Value = 6.539599999999999
$ws.Range("B14").Value = 8.525600000000008
$ws.Range("D14").Value = -8.674600000000003
$ws.Range("C19").Value = -12.68529999999999
$ws.Range("D19").Value = -8.165399999999998
$ws.Range("C21").Value = -12.97
$ws.Range("D21").Value = -9.194799999999987
$ws.Range("B22").Value = 4.783200000000004
$ws.Range("C24").Value = -11.8219
